$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New PSSM values for rows 2-21 (amino acids), columns B-K (positions 0-9)
# Generated from the commit diff: "updates with supplemental figures"
$data = @(
    @(-17.5319753533212,1.725764666851747,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212),
    @(-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212,1.786119457275393,-17.5319753533212,-17.5319753533212),
    @(-17.5319753533212,2.098649059014966,1.977384077298559,-17.5319753533212,3.676030922003461,-17.5319753533212,1.97157285094528,-17.5319753533212,-17.5319753533212,-17.5319753533212),
    @(-17.5319753533212,1.835415934145666,-17.5319753533212,-17.5319753533212,-17.5319753533212,2.255074293590282,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212),
    @(-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212),
    @(2.635493021337801,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212),
    @(-17.5319753533212,-17.5319753533212,-17.5319753533212,2.020462976775399,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212),
    @(3.785135246707882,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212),
    @(-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212,0.5264656614193748,-17.5319753533212,1.948209364911713),
    @(-17.5319753533212,-17.5319753533212,-17.5319753533212,2.84137518996667,-17.5319753533212,3.620666617381592,-17.5319753533212,-17.5319753533212,-17.5319753533212,2.435500006977185),
    @(-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212),
    @(-17.5319753533212,-17.5319753533212,-17.5319753533212,2.523790311284424,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212,1.92337795715255),
    @(-17.5319753533212,-17.5319753533212,1.040402912822723,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212,2.249898190373996),
    @(-17.5319753533212,-17.5319753533212,1.00639927010768,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212),
    @(-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212),
    @(-17.5319753533212,2.408310998014676,2.885783439038977,-17.5319753533212,-17.5319753533212,-17.5319753533212,1.850692996946451,3.028339753129574,4.321920863048043,-17.5319753533212),
    @(-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212,-17.5319753533212,2.360491866422322,1.991389422368652,-17.5319753533212,-17.5319753533212),
    @(-17.5319753533212,-17.5319753533212,1.330507639426023,-17.5319753533212,-17.5319753533212,-17.5319753533212,1.147783059415444,0.8216582513666666,-17.5319753533212,-17.5319753533212),
    @(-17.5319753533212,0.7923326683969405,1.063824814409492,-17.5319753533212,2.851599223388329,-17.5319753533212,1.334678757508934,0.2730190194024142,-17.5319753533212,1.125555939492316),
    @(-17.5319753533212,0.848117854751704,-17.5319753533212,1.596857254294011,-17.5319753533212,1.548751107949206,1.377323007009149,-17.5319753533212,-17.5319753533212,-17.5319753533212)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $col = 2 + $j
        $ws.Cells.Item($row, $col).Value = $rowValues[$j]
    }
}
